# Scentaur Report proof-reading edits
#
# Strategy: for every location that needs a run-level structural change
# (splitting runs, adding/removing <w:proofErr/>, <w:bookmarkStart/End/>,
# <w:lastRenderedPageBreak/>, etc.) we:
#   1. Find.Execute the distinctive plain-text for that location to get a
#      Range over exactly the text we want to replace.
#   2. Clear the text of that range.
#   3. Insert a hand-built run-level OOXML fragment (via Range.InsertXML)
#      at the now-collapsed point, reproducing the desired <w:r>/<w:proofErr/>
#      /<w:bookmarkStart/>/<w:lastRenderedPageBreak/> structure exactly.

function Replace-WithRunXml {
    param($doc, $findText, $xmlBody)

    $sel = $doc.Content
    $found = $sel.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $findText"
    }
    $r = $sel.Duplicate
    $r.Text = ""
    $r2 = $doc.Range($r.Start, $r.Start)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $xmlBody + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r2.InsertXML($xml)
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the _GoBack bookmark after "1.4 " (before "Typical User
#    Experience"). Reconstruct both runs identically minus the bookmark.
# ---------------------------------------------------------------------
Replace-WithRunXml $d "1.4 Typical User Experience" (
    '<w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">1.4 </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Typical User Experience</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 2. "Scentaur is planned to become a web application..." -> split out
#    "The plan is for " + spell-checked "Scentaur" + rest of sentence.
# ---------------------------------------------------------------------
Replace-WithRunXml $d "Scentaur is planned to become a web application through Spring. Therefore, users are enabled to either submit a zip folder or java files directly to a web server. If a Zip folder is submitted, its contents will be extracted to a directory. Otherwise, files will be placed into a directory." (
    '<w:r><w:t xml:space="preserve">The plan is for </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Scentaur</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> to become a web application through Spring. Therefore, users are enabled to either submit a zip folder or java files directly to a web server. If a Zip folder is submitted, its contents will be extracted to a directory. Otherwise, files will be placed into a directory.</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 3. "Once a user is done with Scentaur, the contents of the directory
#    is wiped." -> "...directory " + "are" + " wiped."
# ---------------------------------------------------------------------
Replace-WithRunXml $d "Once a user is done with Scentaur, the contents of the directory is wiped." (
    '<w:r><w:t xml:space="preserve">Once a user is done with Scentaur, the contents of the directory </w:t></w:r>' +
    '<w:r><w:t>are</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> wiped.</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 4. "The constructor will call configureSymbolSolver and parse all
#    source file based on..." -> "...source file" + "s" + " based on..."
# ---------------------------------------------------------------------
Replace-WithRunXml $d "The constructor will call configureSymbolSolver and parse all source file based on JavaParser-JUG-Milano slides." (
    '<w:r><w:t>The constructor will call configureSymbolSolver and parse all source file</w:t></w:r>' +
    '<w:r><w:t>s</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> based on JavaParser-JUG-Milano slides.</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 5/6. Move <w:lastRenderedPageBreak/> from the start of the
#    "Bloater, Abuser, Coupler, Dispensable" run to the start of the
#    "These categories include:" run.
# ---------------------------------------------------------------------
Replace-WithRunXml $d "These categories include:" (
    '<w:r><w:lastRenderedPageBreak/><w:t>These categories include:</w:t></w:r>'
)
Replace-WithRunXml $d "Bloater, Abuser, Coupler, Dispensable" (
    '<w:r><w:t>Bloater, Abuser, Coupler, Dispensable</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 7. ", Coupleable and Dispensable. " + "To ensure..." -> split further
#    with proofErr around "Coupleable" and rewording to "; to ensure...".
# ---------------------------------------------------------------------
Replace-WithRunXml $d ", Coupleable and Dispensable. To ensure that we enable plug-in-play system for the smells within the categories." (
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>Coupleable</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> and Dispensable</w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>;</w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>t</w:t></w:r>' +
    '<w:r><w:t>o ensure that we enable plug-in-play system for the smells within the categories.</w:t></w:r>'
)

Write-Output "Basic text edits done"
